$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9644012944983819
$ws.Range("C2").Value = 0.8608414239482201

$ws.Range("B3").Value = 0.9579288025889967
$ws.Range("C3").Value = 0.8705501618122977

$ws.Range("B4").Value = 0.970873786407767
$ws.Range("C4").Value = 0.8608414239482201

$ws.Range("B5").Value = 0.9644012944983819
$ws.Range("C5").Value = 0.8673139158576052

$ws.Range("B6").Value = 0.9644012944983819
$ws.Range("C6").Value = 0.8543689320388349
